# Insert a new data row at row 815 (2026/02/15, 日, 16, 201), pushing the
# existing rows 815-856 down to 816-857. The workbook grows from
# A1:D856 to A1:D857.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 815:856 down by one row.
$ws.Rows.Item(815).Insert()

# Populate the newly inserted row. Column A holds a date-like string that
# must stay plain text (matching the rest of the sheet, which stores dates
# as literal strings, not Excel date serials). Format it as Text first,
# then strip the formatting afterwards so the resulting cell has no
# explicit style - same as its neighbours.
$ws.Range("A815").NumberFormat = "@"
$ws.Range("A815").Value = "2026/02/15"
$ws.Range("A815").ClearFormats()

$ws.Range("B815").Value = "日"
$ws.Range("C815").Value = 16
$ws.Range("D815").Value = 201
